$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.919.04"
$ws.Range("E2").Value = "  -0.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.630.66"
$ws.Range("E3").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.91"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.33%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Solana
$ws.Range("D8").Value = "23.25"
$ws.Range("E8").Value = "  -0.64%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.02%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.04%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.23%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.861.73"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.619.27"
$ws.Range("E13").Value = "  -0.72%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.87%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.29%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.82"
$ws.Range("E16").Value = "  -1.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.926.67"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "227.81"
$ws.Range("E18").Value = "  -1.22%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -0.52%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.82%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.09%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.07%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.98"
$ws.Range("E23").Value = "  -3.65%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +1.78%  "

# Row 25 - Monero
$ws.Range("D25").Value = "154.68"
$ws.Range("E25").Value = "  -0.18%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.40"
$ws.Range("E29").Value = "  -1.05%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.44%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.37%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  -0.05%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.420.50"
$ws.Range("E33").Value = "  +1.33%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +1.05%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +2.66%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -1.40%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -1.14%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.90%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -0.48%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.70%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -2.29%  "

# Row 42 - Aave
$ws.Range("D42").Value = "65.80"
$ws.Range("E42").Value = "  -1.13%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -1.43%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -0.89%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.770.96"
$ws.Range("E45").Value = "  -0.10%  "

# Row 46 - MXToken
$ws.Range("E46").Value = "  -3.77%  "

# Row 47 - Quant
$ws.Range("D47").Value = "88.61"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +0.44%  "

# Row 49 - was BabyDogeCoin, now Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  -0.39%  "

# Row 50 - was Cronos, now EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  +0.64%  "

# Row 51 - was EnergySwap, now USDD
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.13%  "
